# The diff shows the "Status" column toggling for two rows that share the
# same owner (Matt Medeiros): the "ZZZZAB Controls & Technology" org becomes
# InActive, while the "610 Investments" org becomes Active. (The resulting
# shared-string re-ordering and minor default column width drift in the
# original diff are incidental save artifacts of the source tool, not
# separate edits.)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 24: ZZZZAB Controls & Technology / Matt Medeiros -> Status: Active -> InActive
$ws.Range("C24").Value = "InActive"

# Row 25: 610 Investments / Matt Medeiros -> Status: InActive -> Active
$ws.Range("C25").Value = "Active"
